# Generate Report for Handoff
#
# Updates the status of file "968e5e94-e23d-45b5-ac89-bda0c44d0223.md"
# (row 3 in each sheet) from "Handed back: in sync with en-US" to
# "Ready for handoff", refreshing the related handoff timestamps.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# --- Overview sheet ---
# Columns: A=File Name, B=zh-cn, C=de-de, D=Latest Handoff Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus
$wsOverview.Range("D3").Value = "2016-03-24 06:53:53"

# --- zh-cn sheet ---
# Columns: A=Source File Name, B=File Extension, C=Status,
#          D=Latest Handoff File, E=Latest Handoff Datetime, ...
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("E3").Value = "2016-03-24 06:53:48"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("E3").Value = "2016-03-24 06:53:53"
